$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 17.60633566666667
$ws.Range("H2").Value = 52.819007
$ws.Range("I2").Value = 0.01967183396478544
$ws.Range("J2").Value = 0.01967183396478544
$ws.Range("M2").Value = 17.10933733333333
$ws.Range("N2").Value = 51.328012
$ws.Range("O2").Value = 0.3554368716515803
$ws.Range("P2").Value = 0.3554368716515803
$ws.Range("Q2").Value = 301.2327361248982
$ws.Range("R2").Value = 2711.094625124084
$ws.Range("S2").Value = 0.006992095124092642
$ws.Range("T2").Value = 0.006992095124092639
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 17.60633566666667
$ws.Range("H3").Value = 52.819007
$ws.Range("I3").Value = 0.01967183396478544
$ws.Range("J3").Value = 0.01967183396478544
$ws.Range("O3").Value = 0.2270123898818874
$ws.Range("P3").Value = 0.2270123898818874
$ws.Range("Q3").Value = 192.3929923775792
$ws.Range("R3").Value = 1731.536931398213
$ws.Range("S3").Value = 0.004465750041705629
$ws.Range("T3").Value = 0.004465750041705626
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 17.60633566666667
$ws.Range("H4").Value = 52.819007
$ws.Range("I4").Value = 0.01967183396478544
$ws.Range("J4").Value = 0.01967183396478544
$ws.Range("M4").Value = 11.616679
$ws.Range("N4").Value = 34.850037
$ws.Range("O4").Value = 0.2413299803667016
$ws.Range("P4").Value = 0.2413299803667016
$ws.Range("Q4").Value = 204.5271498059176
$ws.Range("R4").Value = 1840.744348253259
$ws.Range("S4").Value = 0.004747403304498686
$ws.Range("T4").Value = 0.004747403304498684
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 17.60633566666667
$ws.Range("H5").Value = 52.819007
$ws.Range("I5").Value = 0.01967183396478544
$ws.Range("J5").Value = 0.01967183396478544
$ws.Range("M5").Value = 8.482576333333332
$ws.Range("N5").Value = 25.447729
$ws.Range("O5").Value = 0.1762207580998305
$ws.Range("P5").Value = 0.1762207580998305
$ws.Range("Q5").Value = 149.3470862427892
$ws.Range("R5").Value = 1344.123776185103
$ws.Range("S5").Value = 0.003466585494488485
$ws.Range("T5").Value = 0.003466585494488484
$ws.Range("I6").Value = 0.6688940299055509
$ws.Range("J6").Value = 0.6688940299055508
$ws.Range("M6").Value = 17.10933733333333
$ws.Range("N6").Value = 51.328012
$ws.Range("O6").Value = 0.3554368716515803
$ws.Range("P6").Value = 0.3554368716515803
$ws.Range("Q6").Value = 10242.70432369198
$ws.Range("R6").Value = 92184.33891322782
$ws.Range("S6").Value = 0.2377496014560476
$ws.Range("T6").Value = 0.2377496014560476
$ws.Range("I7").Value = 0.6688940299055509
$ws.Range("J7").Value = 0.6688940299055508
$ws.Range("O7").Value = 0.2270123898818874
$ws.Range("P7").Value = 0.2270123898818874
$ws.Range("S7").Value = 0.1518472323065858
$ws.Range("T7").Value = 0.1518472323065858
$ws.Range("I8").Value = 0.6688940299055509
$ws.Range("J8").Value = 0.6688940299055508
$ws.Range("M8").Value = 11.616679
$ws.Range("N8").Value = 34.850037
$ws.Range("O8").Value = 0.2413299803667016
$ws.Range("P8").Value = 0.2413299803667016
$ws.Range("Q8").Value = 6954.460357060498
$ws.Range("R8").Value = 62590.14321354448
$ws.Range("S8").Value = 0.1614241831045105
$ws.Range("T8").Value = 0.1614241831045105
$ws.Range("I9").Value = 0.6688940299055509
$ws.Range("J9").Value = 0.6688940299055508
$ws.Range("M9").Value = 8.482576333333332
$ws.Range("N9").Value = 25.447729
$ws.Range("O9").Value = 0.1762207580998305
$ws.Range("P9").Value = 0.1762207580998305
$ws.Range("Q9").Value = 5078.193245755198
$ws.Range("R9").Value = 45703.73921179679
$ws.Range("S9").Value = 0.1178730130384069
$ws.Range("T9").Value = 0.1178730130384069
$ws.Range("G10").Value = 274.6625416666666
$ws.Range("H10").Value = 823.987625
$ws.Range("I10").Value = 0.3068847498029997
$ws.Range("J10").Value = 0.3068847498029996
$ws.Range("M10").Value = 17.10933733333333
$ws.Range("N10").Value = 51.328012
$ws.Range("O10").Value = 0.3554368716515803
$ws.Range("P10").Value = 0.3554368716515803
$ws.Range("Q10").Value = 4699.294078205721
$ws.Range("R10").Value = 42293.6467038515
$ws.Range("S10").Value = 0.1090781554275561
$ws.Range("T10").Value = 0.1090781554275561
$ws.Range("G11").Value = 274.6625416666666
$ws.Range("H11").Value = 823.987625
$ws.Range("I11").Value = 0.3068847498029997
$ws.Range("J11").Value = 0.3068847498029996
$ws.Range("O11").Value = 0.2270123898818874
$ws.Range("P11").Value = 0.2270123898818874
$ws.Range("Q11").Value = 3001.371170341097
$ws.Range("R11").Value = 27012.34053306988
$ws.Range("S11").Value = 0.06966664047108403
$ws.Range("T11").Value = 0.06966664047108402
$ws.Range("G12").Value = 274.6625416666666
$ws.Range("H12").Value = 823.987625
$ws.Range("I12").Value = 0.3068847498029997
$ws.Range("J12").Value = 0.3068847498029996
$ws.Range("M12").Value = 11.616679
$ws.Range("N12").Value = 34.850037
$ws.Range("O12").Value = 0.2413299803667016
$ws.Range("P12").Value = 0.2413299803667016
$ws.Range("Q12").Value = 3190.666579865791
$ws.Range("R12").Value = 28715.99921879212
$ws.Range("S12").Value = 0.07406049064479805
$ws.Range("T12").Value = 0.07406049064479804
$ws.Range("G13").Value = 274.6625416666666
$ws.Range("H13").Value = 823.987625
$ws.Range("I13").Value = 0.3068847498029997
$ws.Range("J13").Value = 0.3068847498029996
$ws.Range("M13").Value = 8.482576333333332
$ws.Range("N13").Value = 25.447729
$ws.Range("O13").Value = 0.1762207580998305
$ws.Range("P13").Value = 0.1762207580998305
$ws.Range("Q13").Value = 2329.845975594847
$ws.Range("R13").Value = 20968.61378035362
$ws.Range("S13").Value = 0.05407946325956142
$ws.Range("T13").Value = 0.05407946325956142
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 4.071711
$ws.Range("H14").Value = 12.215133
$ws.Range("I14").Value = 0.004549386326664026
$ws.Range("J14").Value = 0.004549386326664025
$ws.Range("M14").Value = 17.10933733333333
$ws.Range("N14").Value = 51.328012
$ws.Range("O14").Value = 0.3554368716515803
$ws.Range("P14").Value = 0.3554368716515803
$ws.Range("Q14").Value = 69.66427702284399
$ws.Range("R14").Value = 626.978493205596
$ws.Range("S14").Value = 0.001617019643883936
$ws.Range("T14").Value = 0.001617019643883935
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 4.071711
$ws.Range("H15").Value = 12.215133
$ws.Range("I15").Value = 0.004549386326664026
$ws.Range("J15").Value = 0.004549386326664025
$ws.Range("O15").Value = 0.2270123898818874
$ws.Range("P15").Value = 0.2270123898818874
$ws.Range("Q15").Value = 44.493566305783
$ws.Range("R15").Value = 400.442096752047
$ws.Range("S15").Value = 0.001032767062511981
$ws.Range("T15").Value = 0.001032767062511981
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 4.071711
$ws.Range("H16").Value = 12.215133
$ws.Range("I16").Value = 0.004549386326664026
$ws.Range("J16").Value = 0.004549386326664025
$ws.Range("M16").Value = 11.616679
$ws.Range("N16").Value = 34.850037
$ws.Range("O16").Value = 0.2413299803667016
$ws.Range("P16").Value = 0.2413299803667016
$ws.Range("Q16").Value = 47.299759667769
$ws.Range("R16").Value = 425.697837009921
$ws.Range("S16").Value = 0.00109790331289437
$ws.Range("T16").Value = 0.00109790331289437
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 4.071711
$ws.Range("H17").Value = 12.215133
$ws.Range("I17").Value = 0.004549386326664026
$ws.Range("J17").Value = 0.004549386326664025
$ws.Range("M17").Value = 8.482576333333332
$ws.Range("N17").Value = 25.447729
$ws.Range("O17").Value = 0.1762207580998305
$ws.Range("P17").Value = 0.1762207580998305
$ws.Range("Q17").Value = 34.53859936477299
$ws.Range("R17").Value = 310.847394282957
$ws.Range("S17").Value = 0.0008016963073737378
$ws.Range("T17").Value = 0.0008016963073737378
